$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings (e.g. "1.00", "2.00") stay as text,
# matching the source data which stores prices as literal strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.339.06"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.525.29"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.53"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.59"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.49"
$ws.Range("E10").Value = "  +4.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.53"
$ws.Range("E11").Value = "  +13.37%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.919.34"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.523.54"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.856"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.160.19"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.48"
$ws.Range("E19").Value = "  +6.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.66"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.98"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.83"
$ws.Range("E24").Value = "  +8.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.07"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  +3.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.83"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.71"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.76"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0789"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.00"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.01"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.01"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.14"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0301"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.018.27"
$ws.Range("E45").Value = "  +2.10%  "
$ws.Range("E46").Value = "  +5.50%  "
$ws.Range("E47").Value = "  +8.60%  "
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.15"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.29"
$ws.Range("E51").Value = "  +3.45%  "
